$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Insert a new "Meta description" paragraph right after the
#    Heading 1 title paragraph at the top of the document.
# ------------------------------------------------------------------
$titlePara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd() -eq "Play 3 Stars Slot for Free - Review and Bonuses" `
            -and $p.Style.NameLocal -eq "Heading 1") {
        $titlePara = $p
        break
    }
}

$labelText = "Meta description"
$restText = ": Experience the rich graphics and progressive jackpots of 3 Stars slot. Read our review and claim free spins, multipliers, and Bonus game rewards."

$titlePara.Range.InsertParagraphAfter()
$metaPara = $titlePara.Next()
$metaPara.Style = "Normal"
$metaPara.Range.Text = $labelText + $restText

$metaStart = $metaPara.Range.Start
$labelRange = $d.Range($metaStart, $metaStart + $labelText.Length)
$labelRange.Bold = 1

# ------------------------------------------------------------------
# 2) Remove the duplicated bold "Play 3 Stars Slot for Free - Review
#    and Bonuses" paragraph that sits near the end of the document
#    (the Normal-styled one, not the Heading 1 title).
# ------------------------------------------------------------------
$dupPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd() -eq "Play 3 Stars Slot for Free - Review and Bonuses" `
            -and $p.Style.NameLocal -eq "Normal") {
        $dupPara = $p
        break
    }
}
$dupPara.Range.Delete()

# ------------------------------------------------------------------
# 3) Replace the text of the closing italic paragraph (formerly the
#    meta-description blurb) with the AI feature-image prompt, while
#    keeping its italic formatting intact.
# ------------------------------------------------------------------
$oldBlurb = "Experience the rich graphics and progressive jackpots of 3 Stars slot. Read our review and claim free spins, multipliers, and Bonus game rewards."
$newBlurb = 'Create a feature image for "3 Stars" that showcases a happy Maya warrior with glasses in a cartoon style. Use bright colors to make the image stand out and include elements of Chinese culture to depict the theme of the game. The Maya warrior should be holding a stack of gold coins, surrounded by Chinese lanterns, and standing in front of a temple. The title of the game, "3 Stars", should be prominently displayed in the image, along with the logo of La JVL, the game development company.'

$blurbPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd() -eq $oldBlurb) {
        $blurbPara = $p
        break
    }
}
$bStart = $blurbPara.Range.Start
$bEnd = $blurbPara.Range.End
$blurbTextRange = $d.Range($bStart, $bEnd - 1)
$blurbTextRange.Text = $newBlurb

Write-Output "done"
